$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old index column (A) entirely, shifting B:F left to A:E.
$ws.Range("A1").EntireColumn.Delete()
